# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-27 11:13:52
#
# Applies the updated attendance pull for the Y2 B25/26 GIT & Liver session
# analysis sheet: reshuffled "Recorded By" name ordering (new pull order from
# source system), two sessions (BIOCHEMISTRY LAB/CBL #2 and PATHOLOGY
# LAB/MUSEUM #1 on 27/11) flipping from Pending/Not Recorded to Recorded with
# their new attendee + headcount, and the downstream session/coverage
# statistics recalculated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Recorded By" reshuffles (same people, new order from the refreshed pull)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System"
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# BIOCHEMISTRY LAB/CBL, session 2 (row 8): Not Recorded -> Recorded
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Range("G8").Value = "AbeerRagheb@med.asu.edu.eg"
$ws.Range("H8").Value = "1/251"
$ws.Range("I8").Value = "Recorded"

# ---------------------------------------------------------------------------
# PATHOLOGY LAB/MUSEUM, session 1 on 27/11 (row 23): Pending -> Recorded
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A23:I23").PasteSpecial(-4122)
$ws.Range("G23").Value = "menna-allah.gamil@med.asu.edu.eg"
$ws.Range("H23").Value = "71/251"
$ws.Range("I23").Value = "Recorded"

# ---------------------------------------------------------------------------
# Class Statistics box (ANATOMY) - Recorded/Missing/Pending session counts
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 12

# L9/L10 hold literal percentage text (not numeric %), so force text entry
# via a quote-prefix and then restore the plain "Value" cell format so no
# stray number format sticks to the cell.
$ws.Range("L9").Value = "'55.2%"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").Value = "'23.4%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Group Statistics box (row 15) mirrors the same recount
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 16
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 12

$ws.Range("R15").Value = "'55.2%"
$ws.Range("K15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").Value = "'23.4%"
$ws.Range("K15").Copy()
$ws.Range("S15").PasteSpecial(-4122)
